$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing the existing data rows (2-25) down to (3-26)
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

function Set-TextCell($rng, $val) {
    # Force the value to be stored as text (not auto-converted to a date serial)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Populate the newly inserted row 2 with the new RPA dataset record (2023-10-19, 신성에스티)
Set-TextCell $ws.Range("A2") "2023-10-19"
$ws.Range("B2").Value = "신성에스티"
$ws.Range("C2").Value = "코스닥"
$ws.Range("D2").Value = 520
$ws.Range("E2").Value = "미래"
$ws.Range("F2").Value = 520
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "대표"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 26000
$ws.Range("N2").Value = 100
Set-TextCell $ws.Range("O2") "2023-10-10"
Set-TextCell $ws.Range("P2") "2023-10-13"
$ws.Range("Q2").Value = 1300000
